$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "test461"
$ws.Range("C2").Value = "narendra695"
$ws.Range("D2").Value = "m4s$!W5E"
$ws.Range("B2").Value = 23071214
